$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Touch B1 (without altering content) to keep the sheet dimension starting at row 1
$ws.Range("B1").Font.Bold = $ws.Range("B1").Font.Bold

# Update D10: step 4 text changed
$ws.Range("D10").Value = "4. Obtém informação"

# Add new step 5 in D11
$ws.Range("D11").Value = "5. Mostra lista de carros comprados"

# Add post-condition text in C5
$ws.Range("C5").Value = "Verificou lista de carros comprados"

# Update the sheet view: top-left cell and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D10").Select()
